$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.761.51'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.306.64'
$ws.Range('E3').Value = '  +3.40%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '269.94'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = '92.86'
$ws.Range('E6').Value = '  +4.89%  '
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').Value = '44.68'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '8.08'
$ws.Range('E12').Value = '  +6.80%  '
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '2.656.54'
$ws.Range('E14').Value = '  +3.73%  '
$ws.Range('D15').Value = '15.30'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = '0.848'
$ws.Range('E16').Value = '  +6.86%  '
$ws.Range('D17').Value = '2.209.53'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '43.744.46'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').Value = '6.24'
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('D21').Value = '71.33'
$ws.Range('E21').Value = '  +1.73%  '
$ws.Range('D22').Value = '239.16'
$ws.Range('E22').Value = '  +2.83%  '
$ws.Range('E23').Value = '  -4.32%  '
$ws.Range('D24').Value = '9.74'
$ws.Range('E24').Value = '  +9.44%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '11.29'
$ws.Range('E26').Value = '  +3.90%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = '2.49'
$ws.Range('E27').Value = '  -4.28%  '
$ws.Range('D28').Value = '2.34'
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('E29').Value = '  -4.76%  '
$ws.Range('D30').Value = '39.00'
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').Value = '22.44'
$ws.Range('E31').Value = '  +8.74%  '
$ws.Range('D32').Value = '171.65'
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('D33').Value = '0.0893'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('D34').Value = '5.60'
$ws.Range('E34').Value = '  +2.99%  '
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = '4.48'
$ws.Range('E37').Value = '  +2.14%  '
$ws.Range('D38').Value = '0.0350'
$ws.Range('E38').Value = '  -2.23%  '
$ws.Range('D39').Value = '3.41'
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('E40').Value = '  +15.36%  '
$ws.Range('D41').Value = '2.30'
$ws.Range('E41').Value = '  +9.06%  '
$ws.Range('D42').Value = '12.13'
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = '1.29'
$ws.Range('E43').Value = '  +14.19%  '
$ws.Range('D44').Value = '5.43'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').Value = '61.30'
$ws.Range('E45').Value = '  -6.54%  '
$ws.Range('D46').Value = '8.87'
$ws.Range('E46').Value = '  +6.68%  '
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('D48').Value = '100.01'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = '2.534.35'
$ws.Range('E50').Value = '  +3.70%  '
$ws.Range('D51').Value = '0.423'
$ws.Range('E51').Value = '  -3.86%  '
